# Weekly price update: insert a new record at the top of the data block
# (row 583), pushing the existing historical rows down by one. This
# mirrors how the source feed prepends the latest week's observation to
# the existing "Pepino dulce" time series held in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 583:706 down to 584:707, leaving a blank row 583 behind
# (this also grows the sheet's used range from R706 to R707, matching
# the updated <dimension> in the workbook).
$ws.Rows.Item(583).Insert()

# Populate the newly inserted row with this week's observation.
$newRow = 583

$ws.Cells.Item($newRow, 1).Value = 8
$ws.Cells.Item($newRow, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($newRow, 3).Value = "Coquimbo"
$ws.Cells.Item($newRow, 4).Value = 45244
$ws.Cells.Item($newRow, 5).Value = 4
$ws.Cells.Item($newRow, 6).Value = 100112043
$ws.Cells.Item($newRow, 7).Value = "Pepino dulce"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Segunda"
$ws.Cells.Item($newRow, 10).Value = 200
$ws.Cells.Item($newRow, 11).Value = 19000
$ws.Cells.Item($newRow, 12).Value = 20000
$ws.Cells.Item($newRow, 13).Value = 19500
$ws.Cells.Item($newRow, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item($newRow, 15).Value = "Provincia de Limarí"
$ws.Cells.Item($newRow, 16).Value = 1083
$ws.Cells.Item($newRow, 17).Value = 18
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
